# "10 years Finalization data"
#
# The monthly station table (header + 32 daily rows, A9:K40) living inside
# the single worksheet "Data Harian - Table" gets split out onto its own
# worksheet ("Sheet1"), which becomes the new active/selected tab. The
# original sheet keeps all of its rows untouched but the on-screen
# selection/scroll position moves down to frame the data table, and the
# now-unused logo picture is removed.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- remove the station logo picture from the original sheet ---------
if ($ws1.Shapes.Count -gt 0) {
    for ($i = $ws1.Shapes.Count; $i -ge 1; $i--) {
        $ws1.Shapes.Item($i).Delete()
    }
}

# --- add the new worksheet right after the original one ---------------
$ws2 = $wb.Worksheets.Add($null, $ws1)

# --- copy the data table (A9:K40) into the new sheet at A1:K32 --------
$src = $ws1.Range("A9:K40")

# formats first (keeps the existing header/body cell style indices)...
$src.Copy()
$ws2.Range("A1").PasteSpecial(-4122)

# ...then values + number formats on top.
$src.Copy()
$ws2.Range("A1").PasteSpecial(-4163)

# Matches the source table's taller, word-wrapped body rows.
$ws2.Range("A2:K32").RowHeight = 28.8

# --- view/selection bookkeeping ----------------------------------------
$ws1.Range("A9:K40").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7

$ws2.Range("A1:K32").Select() | Out-Null
$ws2.Activate() | Out-Null
